# This script reorders several paragraph text blocks within the LOM3267 document.
# The document structure (paragraph styles, runs, line breaks, bold/italic runs) stays
# unchanged; only the w:t text content of certain runs moves to different paragraphs.
# Because several of these moves form cycles (A->B->C->A), a two-phase swap through
# unique placeholder tokens is used to avoid a later Find matching text that was only
# just written by an earlier step.

$d = $word.ActiveDocument

# ---- Phase 1: stash original text of every run that changes behind a unique placeholder ----
$r = $d.Content
$r.Find.Execute("O Trabalho de Graduação (TG) tem por objetivo a integração, o aprofundamento e aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de tarefas que fazem parte do perfil de atuação profissional do engenheiro físico.", $true) | Out-Null
$r.Text = "@@PH0@@"

$r = $d.Content
$r.Find.Execute("The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to perform tasks that are part of the professional performance profile of the physical engineer.", $true) | Out-Null
$r.Text = "@@PH1@@"

$r = $d.Content
$r.Find.Execute("5840730 - Antonio Jefferson da Silva Machado", $true) | Out-Null
$r.Text = "@@PH2@@"

$r = $d.Content
$r.Find.Execute("1176388 - Luiz Tadeu Fernandes Eleno", $true) | Out-Null
$r.Text = "@@PH3@@"

$r = $d.Content
$r.Find.Execute("Elaboração, com a orientação de um professor supervisor, de uma proposta de projeto em tema ligado à área de ciência e tecnologia.", $true) | Out-Null
$r.Text = "@@PH4@@"

$r = $d.Content
$r.Find.Execute("Preparation, with the guidance of a supervising professor, of a project proposal on a topic related to the area of science and technology.", $true) | Out-Null
$r.Text = "@@PH5@@"

$r = $d.Content
$r.Find.Execute("O aluno deve procurar um professor ou profissional com formação na área de engenharia ou áreas correlatas, para a elaboração de uma proposta de projeto contendo motivação e objetivos, fundamentação teórica e cronograma de execução. O projeto propriamente dito será desenvolvido e defendido na disciplina Trabalho de Graduação II.", $true) | Out-Null
$r.Text = "@@PH6@@"

$r = $d.Content
$r.Find.Execute("O aluno deve apresentar a proposta de trabalho à uma banca formada pelo responsável pela disciplina e professores ou profissionais da área.", $true) | Out-Null
$r.Text = "@@PH7@@"

$r = $d.Content
$r.Find.Execute("Avaliação e atribuição de nota do Trabalho de Graduação por uma comissão de professores.", $true) | Out-Null
$r.Text = "@@PH8@@"

$r = $d.Content
$r.Find.Execute("A critério da banca de avaliação poderá ser estabelecido um prazo para readequação e reapresentação do plano de trabalho.", $true) | Out-Null
$r.Text = "@@PH9@@"

$r = $d.Content
$r.Find.Execute("A ser definida no plano de trabalho.", $true) | Out-Null
$r.Text = "@@PH10@@"

# ---- Phase 2: write the final text for each placeholder ----
$r = $d.Content
$r.Find.Execute("@@PH0@@", $true) | Out-Null
$r.Text = "Elaboração, com a orientação de um professor supervisor, de uma proposta de projeto em tema ligado à área de ciência e tecnologia."

$r = $d.Content
$r.Find.Execute("@@PH1@@", $true) | Out-Null
$r.Text = "Preparation, with the guidance of a supervising professor, of a project proposal on a topic related to the area of science and technology."

$r = $d.Content
$r.Find.Execute("@@PH2@@", $true) | Out-Null
$r.Text = "O Trabalho de Graduação (TG) tem por objetivo a integração, o aprofundamento e aplicação dos conhecimentos adquiridos ao longo do curso, preparando e desenvolvendo a capacidade do aluno para a realização de tarefas que fazem parte do perfil de atuação profissional do engenheiro físico."

$r = $d.Content
$r.Find.Execute("@@PH3@@", $true) | Out-Null
$r.Text = "O aluno deve procurar um professor ou profissional com formação na área de engenharia ou áreas correlatas, para a elaboração de uma proposta de projeto contendo motivação e objetivos, fundamentação teórica e cronograma de execução. O projeto propriamente dito será desenvolvido e defendido na disciplina Trabalho de Graduação II."

$r = $d.Content
$r.Find.Execute("@@PH4@@", $true) | Out-Null
$r.Text = "O aluno deve apresentar a proposta de trabalho à uma banca formada pelo responsável pela disciplina e professores ou profissionais da área."

$r = $d.Content
$r.Find.Execute("@@PH5@@", $true) | Out-Null
$r.Text = "The Graduation Work (TG) aims to integrate, deepen and apply the knowledge acquired throughout the course, preparing and developing the student's ability to perform tasks that are part of the professional performance profile of the physical engineer."

$r = $d.Content
$r.Find.Execute("@@PH6@@", $true) | Out-Null
$r.Text = "Avaliação e atribuição de nota do Trabalho de Graduação por uma comissão de professores."

$r = $d.Content
$r.Find.Execute("@@PH7@@", $true) | Out-Null
$r.Text = "A critério da banca de avaliação poderá ser estabelecido um prazo para readequação e reapresentação do plano de trabalho."

$r = $d.Content
$r.Find.Execute("@@PH8@@", $true) | Out-Null
$r.Text = "A ser definida no plano de trabalho."

$r = $d.Content
$r.Find.Execute("@@PH9@@", $true) | Out-Null
$r.Text = "5840730 - Antonio Jefferson da Silva Machado"

$r = $d.Content
$r.Find.Execute("@@PH10@@", $true) | Out-Null
$r.Text = "1176388 - Luiz Tadeu Fernandes Eleno"

